$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This script reflects a "Generate Report for Handoff" status update:
#  - The status of b.md flips from "Handed back: in sync with en-US" to
#    "Ready for handoff" (Overview sheet + both locale sheets).
#  - For zh-cn and de-de, the "b.md" handback details move forward: new
#    handback xlf filenames, a later handback datetime, the Content
#    Duplicate flag flips to False, and a new validation error message is
#    recorded because the handback file version is stale.
# ---------------------------------------------------------------------------

$newErrorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b771da1e7451c083cbec19fa57c8dad9ed095480/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4272a305ba4735f45286d376cb0a28e92f5409c3/e2e/b.md."

# --- Overview sheet: row 3 is the "b.md" row ------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 12:45:53"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 ("a.md") - status label only
$wsZhCn.Range("C2").Value = "Ready for handoff"
# Row 3 ("b.md") - status, content-duplicate flag, handback file/datetime, error detail
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text storage so "False" isn't coerced to a real
# Excel boolean (the column stores True/False as plain localized text); reset
# the style afterwards so the quote-prefix flag doesn't linger on the cell.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 12:45:48"
$wsZhCn.Range("P3").Value = $newErrorMessage
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 ("b.md") - status, content-duplicate flag, handback file/datetime, error detail
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 12:45:53"
$wsDeDe.Range("P3").Value = $newErrorMessage
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
